$d = $word.ActiveDocument

# --- Change 1: remove the first two "Intro" paragraphs ---
$p1 = $d.Paragraphs(1)
$p2 = $d.Paragraphs(2)
$introRange = $d.Range($p1.Range.Start, $p2.Range.End)
$introRange.Delete()

# --- Change 2: merge the trailing "system to another." + quote + citation runs ---
$text = $d.Content.Text
$startIdx = $text.IndexOf("system to another.")
$endIdx = $text.IndexOf("(Microsoft, n.d.)") + "(Microsoft, n.d.)".Length
$mergeRange = $d.Range($startIdx, $endIdx)
$mergeRange.Delete()
$insertAt = $d.Range($startIdx, $startIdx)
$insertAt.InsertAfter("system to another." + [char]0x201D + " (Microsoft, n.d.)")

# --- Change 3: add the four new tool bullet points after the Airbyte bullet ---
$text = $d.Content.Text
$toolIdx = $text.IndexOf(") tool.") + ") tool.".Length
$insertPoint = $d.Range($toolIdx, $toolIdx)
$newParasBody = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Fivetran</w:t></w:r><w:r><w:t>: automated ETL tool</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Matillion</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>cloud-native data integration platform</w:t></w:r><w:r><w:t xml:space="preserve"> that facilitates ETL</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>IBM Informix</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>relational database management system (RDBMS)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>CloudFuze</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>cloud content management platform</w:t></w:r><w:r><w:t xml:space="preserve">. Can </w:t></w:r><w:r><w:t>migrate and govern data across multiple cloud storage providers</w:t></w:r></w:p>'
$newParasXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $newParasBody + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($newParasXml)

# --- Change 4: mark the last picture (the one added without noProof) as NoProofing ---
$lastShape = $d.InlineShapes(7)
$lastShape.Range.NoProofing = $true

# --- Change 5: enable the (latent) Heading 3 style definition, mirroring what Word does
#     the first time a built-in heading style is applied, without leaving stray content ---
$tempPara = $d.Paragraphs.Add()
$tempPara.Range.Text = "zzz_tmp_heading3_zzz"
$tempPara.Range.Style = "Heading 3"
$tempPara.Range.Delete()
